$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 37656
$ws.Cells.Item(2, 4).Value = 54460524
$ws.Cells.Item(3, 3).Value = 90822
$ws.Cells.Item(3, 4).Value = 133139845
$ws.Cells.Item(4, 3).Value = 31128
$ws.Cells.Item(4, 4).Value = 46100954
$ws.Cells.Item(5, 3).Value = 8681
$ws.Cells.Item(5, 4).Value = 12902063
$ws.Cells.Item(6, 3).Value = 1986
$ws.Cells.Item(6, 4).Value = 2951506
$ws.Cells.Item(11, 3).Value = 41221
$ws.Cells.Item(11, 4).Value = 55936285
$ws.Cells.Item(12, 3).Value = 9635
$ws.Cells.Item(12, 4).Value = 13936170
$ws.Cells.Item(13, 3).Value = 25913
$ws.Cells.Item(13, 4).Value = 38005540
$ws.Cells.Item(14, 3).Value = 8304
$ws.Cells.Item(14, 4).Value = 12324218
$ws.Cells.Item(19, 3).Value = 10204
$ws.Cells.Item(19, 4).Value = 13512858
$ws.Cells.Item(20, 3).Value = 13359
$ws.Cells.Item(20, 4).Value = 19291294
$ws.Cells.Item(21, 3).Value = 31606
$ws.Cells.Item(21, 4).Value = 46384039
$ws.Cells.Item(22, 3).Value = 10210
$ws.Cells.Item(22, 4).Value = 15178055
$ws.Cells.Item(23, 3).Value = 2635
$ws.Cells.Item(23, 4).Value = 3917682
$ws.Cells.Item(24, 3).Value = 505
$ws.Cells.Item(24, 4).Value = 751592
$ws.Cells.Item(26, 3).Value = 11665
$ws.Cells.Item(26, 4).Value = 15583911
$ws.Cells.Item(27, 3).Value = 7629
$ws.Cells.Item(27, 4).Value = 11052062
$ws.Cells.Item(28, 3).Value = 22444
$ws.Cells.Item(28, 4).Value = 32944405
$ws.Cells.Item(29, 3).Value = 7797
$ws.Cells.Item(29, 4).Value = 11603133
$ws.Cells.Item(30, 3).Value = 1957
$ws.Cells.Item(30, 4).Value = 2919999
$ws.Cells.Item(33, 3).Value = 8282
$ws.Cells.Item(33, 4).Value = 10942109
$ws.Cells.Item(34, 3).Value = 3235
$ws.Cells.Item(34, 4).Value = 4669637
$ws.Cells.Item(35, 3).Value = 7811
$ws.Cells.Item(35, 4).Value = 11406956
$ws.Cells.Item(36, 3).Value = 3172
$ws.Cells.Item(36, 4).Value = 4700961
$ws.Cells.Item(37, 3).Value = 828
$ws.Cells.Item(37, 4).Value = 1233223
$ws.Cells.Item(40, 3).Value = 2467
$ws.Cells.Item(40, 4).Value = 3335331
$ws.Cells.Item(41, 3).Value = 17200
$ws.Cells.Item(41, 4).Value = 24872665
$ws.Cells.Item(42, 3).Value = 51002
$ws.Cells.Item(42, 4).Value = 74771431
$ws.Cells.Item(43, 3).Value = 18985
$ws.Cells.Item(43, 4).Value = 28199943
$ws.Cells.Item(44, 3).Value = 5599
$ws.Cells.Item(44, 4).Value = 8338478
$ws.Cells.Item(49, 3).Value = 16649
$ws.Cells.Item(49, 4).Value = 22173203
$ws.Cells.Item(51, 3).Value = 6855
$ws.Cells.Item(51, 4).Value = 10077424
$ws.Cells.Item(52, 3).Value = 2337
$ws.Cells.Item(52, 4).Value = 3490418
$ws.Cells.Item(56, 3).Value = 6866
$ws.Cells.Item(56, 4).Value = 9449525
$ws.Cells.Item(57, 3).Value = 932
$ws.Cells.Item(57, 4).Value = 1367579
$ws.Cells.Item(58, 3).Value = 2346
$ws.Cells.Item(58, 4).Value = 3477917
$ws.Cells.Item(59, 3).Value = 936
$ws.Cells.Item(59, 4).Value = 1393501
$ws.Cells.Item(60, 3).Value = 320
$ws.Cells.Item(60, 4).Value = 479758
$ws.Cells.Item(61, 3).Value = 102
$ws.Cells.Item(61, 4).Value = 152850
$ws.Cells.Item(63, 3).Value = 1377
$ws.Cells.Item(63, 4).Value = 1936706
$ws.Cells.Item(64, 3).Value = 15317
$ws.Cells.Item(64, 4).Value = 22125754
$ws.Cells.Item(65, 3).Value = 44598
$ws.Cells.Item(65, 4).Value = 65263308
$ws.Cells.Item(66, 3).Value = 15674
$ws.Cells.Item(66, 4).Value = 23294182
$ws.Cells.Item(67, 3).Value = 4562
$ws.Cells.Item(67, 4).Value = 6795292
$ws.Cells.Item(72, 3).Value = 15054
$ws.Cells.Item(72, 4).Value = 19853457
$ws.Cells.Item(73, 3).Value = 51191
$ws.Cells.Item(73, 4).Value = 74494111
$ws.Cells.Item(74, 3).Value = 145542
$ws.Cells.Item(74, 4).Value = 214419948
$ws.Cells.Item(75, 3).Value = 63448
$ws.Cells.Item(75, 4).Value = 94547770
$ws.Cells.Item(76, 3).Value = 20270
$ws.Cells.Item(76, 4).Value = 30285686
$ws.Cells.Item(77, 3).Value = 4798
$ws.Cells.Item(77, 4).Value = 7168723
$ws.Cells.Item(84, 3).Value = 50636
$ws.Cells.Item(84, 4).Value = 68895541
$ws.Cells.Item(85, 3).Value = 4577
$ws.Cells.Item(85, 4).Value = 6631511
$ws.Cells.Item(86, 3).Value = 11535
$ws.Cells.Item(86, 4).Value = 16946765
$ws.Cells.Item(87, 3).Value = 3875
$ws.Cells.Item(87, 4).Value = 5775083
$ws.Cells.Item(92, 3).Value = 5383
$ws.Cells.Item(92, 4).Value = 7239982
$ws.Cells.Item(94, 3).Value = 5134
$ws.Cells.Item(94, 4).Value = 7561020
$ws.Cells.Item(95, 3).Value = 1936
$ws.Cells.Item(95, 4).Value = 2883937
$ws.Cells.Item(96, 3).Value = 687
$ws.Cells.Item(96, 4).Value = 1029460
$ws.Cells.Item(97, 3).Value = 184
$ws.Cells.Item(97, 4).Value = 275113
$ws.Cells.Item(100, 3).Value = 3526
$ws.Cells.Item(100, 4).Value = 4665599
$ws.Cells.Item(101, 3).Value = 598
$ws.Cells.Item(101, 4).Value = 890664
$ws.Cells.Item(106, 3).Value = 10731
$ws.Cells.Item(106, 4).Value = 15568451
$ws.Cells.Item(107, 3).Value = 29127
$ws.Cells.Item(107, 4).Value = 42797412
$ws.Cells.Item(108, 3).Value = 9760
$ws.Cells.Item(108, 4).Value = 14513650
$ws.Cells.Item(109, 3).Value = 2678
$ws.Cells.Item(109, 4).Value = 3993207
$ws.Cells.Item(110, 3).Value = 489
$ws.Cells.Item(110, 4).Value = 728546
$ws.Cells.Item(113, 3).Value = 9766
$ws.Cells.Item(113, 4).Value = 12903503
$ws.Cells.Item(114, 3).Value = 30349
$ws.Cells.Item(114, 4).Value = 43767040
$ws.Cells.Item(115, 3).Value = 66010
$ws.Cells.Item(115, 4).Value = 96610451
$ws.Cells.Item(116, 3).Value = 21326
$ws.Cells.Item(116, 4).Value = 31693440
$ws.Cells.Item(117, 3).Value = 6050
$ws.Cells.Item(117, 4).Value = 9013561
$ws.Cells.Item(123, 3).Value = 25774
$ws.Cells.Item(123, 4).Value = 34436279
$ws.Cells.Item(124, 3).Value = 35880
$ws.Cells.Item(124, 4).Value = 51787559
$ws.Cells.Item(125, 3).Value = 76593
$ws.Cells.Item(125, 4).Value = 112006996
$ws.Cells.Item(126, 3).Value = 23795
$ws.Cells.Item(126, 4).Value = 35316052
$ws.Cells.Item(127, 3).Value = 6376
$ws.Cells.Item(127, 4).Value = 9475051
$ws.Cells.Item(128, 3).Value = 1229
$ws.Cells.Item(128, 4).Value = 1827911
$ws.Cells.Item(132, 3).Value = 31699
$ws.Cells.Item(132, 4).Value = 42104210
$ws.Cells.Item(133, 3).Value = 13184
$ws.Cells.Item(133, 4).Value = 19082907
$ws.Cells.Item(134, 3).Value = 32244
$ws.Cells.Item(134, 4).Value = 47362472
$ws.Cells.Item(135, 3).Value = 11454
$ws.Cells.Item(135, 4).Value = 17019042
$ws.Cells.Item(140, 3).Value = 10786
$ws.Cells.Item(140, 4).Value = 14384802
$ws.Cells.Item(141, 3).Value = 34929
$ws.Cells.Item(141, 4).Value = 50443593
$ws.Cells.Item(142, 3).Value = 81035
$ws.Cells.Item(142, 4).Value = 118731387
$ws.Cells.Item(143, 3).Value = 24286
$ws.Cells.Item(143, 4).Value = 36083555
$ws.Cells.Item(144, 3).Value = 6370
$ws.Cells.Item(144, 4).Value = 9504567
$ws.Cells.Item(148, 3).Value = 29091
$ws.Cells.Item(148, 4).Value = 39253971
